{"js": "// Resume edits:\n// 1. \"Statistical analysis of\" -> \"Statistical analyses of\"\n//    (typo fix in the HIV/DNA bullet point).\n// 2. \"Software Frameworks: R Markdown, Jupyter Notebooks, LaTeX, and\n//    Node.js.\" -> \"Tools: R Markdown, Jupyter Notebooks, LaTeX, Node.js,\n//    Amazon AWS, Google Cloud Platform, Theano, and TensorFlow.\"\n//    (relabel the bold \"Software Frameworks:\" run to \"Tools:\" and extend\n//    the technical-skills list with additional tools).\n\nconst body = context.document.body;\n\n// --- Change 1: \"analysis\" -> \"analyses\" -------------------------------\nconst analysisHits = body.search(\"Statistical analysis of\", { matchCase: true });\nanalysisHits.load(\"items\");\nawait context.sync();\nif (analysisHits.items.length > 0) {\n  analysisHits.items[0].insertText(\"Statistical analyses of\", \"Replace\");\n  await context.sync();\n}\n\n// --- Change 2: \"Software Frameworks:\" -> \"Tools:\" list ---------------\n// 2a. Re-label the bold run \"Software Frameworks: \" -> \"Tools: \" (keeps\n//     the run's existing bold formatting intact).\nconst labelHits = body.search(\"Software Frameworks: \", { matchCase: true });\nlabelHits.load(\"items\");\nawait context.sync();\nif (labelHits.items.length > 0) {\n  labelHits.items[0].insertText(\"Tools: \", \"Replace\");\n  await context.sync();\n}\n\n// 2b. Drop the \"and \" that precedes \"Node.js\" in the middle of the list\n//     (the list becomes comma-separated with \"and\" only before the final\n//     item now).\nconst andNodeHits = body.search(\"LaTeX, and Node.js\", { matchCase: true });\nandNodeHits.load(\"items\");\nawait context.sync();\nif (andNodeHits.items.length > 0) {\n  andNodeHits.items[0].insertText(\"LaTeX, Node.js\", \"Replace\");\n  await context.sync();\n}\n\n// 2c. Replace the trailing \".\" after \"Node.js\" with the extended list of\n//     additional tools/platforms.\nconst trailingHits = body.search(\"Node.js.\", { matchCase: true });\ntrailingHits.load(\"items\");\nawait context.sync();\nif (trailingHits.items.length > 0) {\n  trailingHits.items[0].insertText(\n    \"Node.js, Amazon AWS, Google Cloud Platform, Theano, and TensorFlow.\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n", "ps1": "# Resume edits:\n# 1. \"Statistical analysis of\" -> \"Statistical analyses of\"\n#    (typo fix in the HIV/DNA bullet point).\n# 2. \"Software Frameworks: R Markdown, Jupyter Notebooks, LaTeX, and\n#    Node.js.\" -> \"Tools: R Markdown, Jupyter Notebooks, LaTeX, Node.js,\n#    Amazon AWS, Google Cloud Platform, Theano, and TensorFlow.\"\n#    (relabel the bold \"Software Frameworks:\" run to \"Tools:\" and extend\n#    the technical-skills list with additional tools).\n\n$d = $word.ActiveDocument\n\nfunction Replace-DocText($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $replaceText\n    # MatchCase=True, Wrap=wdFindStop(0), Replace=wdReplaceAll(2) ; preserves\n    # the source run's formatting (bold label stays bold, body text stays\n    # regular) since no explicit Format is requested.\n    $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 0, $false, $find.Replacement.Text, 2) | Out-Null\n}\n\n# --- Change 1: \"analysis\" -> \"analyses\" --------------------------------\nReplace-DocText \"Statistical analysis of\" \"Statistical analyses of\"\n\n# --- Change 2: \"Software Frameworks:\" -> \"Tools:\" list -----------------\n# 2a. Re-label the bold run \"Software Frameworks: \" -> \"Tools: \".\nReplace-DocText \"Software Frameworks: \" \"Tools: \"\n\n# 2b. Drop the \"and \" that used to precede \"Node.js\" mid-list.\nReplace-DocText \"LaTeX, and Node.js\" \"LaTeX, Node.js\"\n\n# 2c. Replace the trailing \".\" after \"Node.js\" with the extended list of\n#     additional tools/platforms.\nReplace-DocText \"Node.js.\" \"Node.js, Amazon AWS, Google Cloud Platform, Theano, and TensorFlow.\"\n"}
